$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.580.05"
$ws.Range("E2").Value = "  +5.51%  "
$ws.Range("D3").Value = "3.183.86"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("D5").Value = "'402.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.11%  "
$ws.Range("D10").Value = "'39.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.66%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.140"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.0885"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("D13").Value = "3.679.90"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "'19.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").Value = "'8.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").Value = "'1.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.35%  "
$ws.Range("D17").Value = "3.190.23"
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").Value = "'10.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "54.537.75"
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("D20").Value = "'3.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.11%  "
$ws.Range("D21").Value = "'12.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("D22").Value = "'0.0000100"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("D23").Value = "'72.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.83%  "
$ws.Range("D24").Value = "'275.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("E25").Value = "  +4.85%  "
$ws.Range("D26").Value = "'8.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "'27.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").Value = "'7.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("D32").Value = "'11.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.65%  "
$ws.Range("D33").Value = "'0.0503"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.90%  "
$ws.Range("D34").Value = "'36.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.66%  "
$ws.Range("D36").Value = "'51.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +6.90%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  +9.81%  "
$ws.Range("D40").Value = "'4.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.69%  "
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").Value = "'131.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").Value = "'22.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "'2.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "2.091.13"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'0.0347"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.71%  "
$ws.Range("D51").Value = "'0.0519"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.72%  "
